$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates derived from the authoritative diff (cryptos.xlsx refresh).
# Each entry carries the row number plus whichever of B (Coin), C (Link),
# D (Price) and E (Volume 1h) columns changed for that row.
$updates = @(
    @{Row=2; D="42.062.64"; E="  -1.54%  "},
    @{Row=3; D="2.293.46"; E="  -2.16%  "},
    @{Row=4; E="  -0.20%  "},
    @{Row=5; D="312.84"; E="  -1.89%  "},
    @{Row=6; D="105.35"; E="  +0.29%  "},
    @{Row=7; E="  -1.67%  "},
    @{Row=8; E="  -0.16%  "},
    @{Row=9; E="  -1.19%  "},
    @{Row=10; E="  -1.63%  "},
    @{Row=11; D="0.0910"; E="  -1.16%  "},
    @{Row=12; E="  +0.22%  "},
    @{Row=13; E="  +1.19%  "},
    @{Row=14; D="0.971"; E="  -1.46%  "},
    @{Row=15; D="15.33"; E="  -3.42%  "},
    @{Row=16; D="2.641.84"; E="  -2.07%  "},
    @{Row=17; D="2.292.12"; E="  -0.85%  "},
    @{Row=18; D="41.910.67"; E="  -1.85%  "},
    @{Row=19; D="7.60"; E="  -2.04%  "},
    @{Row=20; E="  -0.68%  "},
    @{Row=21; D="72.62"; E="  -6.22%  "},
    @{Row=22; D="3.51"; E="  -1.24%  "},
    @{Row=23; D="258.80"; E="  -0.69%  "},
    @{Row=24; D="2.31"; E="  -0.31%  "},
    @{Row=25; D="9.73"; E="  +0.72%  "},
    @{Row=26; E="  +0.61%  "},
    @{Row=27; D="10.95"; E="  -3.96%  "},
    @{Row=28; E="  +2.44%  "},
    @{Row=29; D="22.69"; E="  -2.10%  "},
    @{Row=30; D="36.09"; E="  -0.43%  "},
    @{Row=31; D="165.14"; E="  -5.67%  "},
    @{Row=32; D="0.0889"; E="  -0.22%  "},
    @{Row=33; E="  -2.80%  "},
    @{Row=34; D="5.87"; E="  -3.45%  "},
    @{Row=35; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.130"; E="  -0.27%  "},
    @{Row=36; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.119"; E="  +6.19%  "},
    @{Row=37; E="  +1.54%  "},
    @{Row=38; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.90"; E="  +10.01%  "},
    @{Row=39; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0350"; E="  -1.28%  "},
    @{Row=40; D="3.60"; E="  -3.91%  "},
    @{Row=41; D="99.04"; E="  +17.05%  "},
    @{Row=42; E="  +1.00%  "},
    @{Row=43; D="70.58"; E="  +0.23%  "},
    @{Row=44; D="0.226"; E="  -2.21%  "},
    @{Row=46; D="12.13"; E="  +2.55%  "},
    @{Row=47; D="112.51"; E="  -1.96%  "},
    @{Row=48; D="78.30"; E="  +6.70%  "},
    @{Row=49; E="  -0.71%  "},
    @{Row=50; D="5.31"; E="  -3.58%  "},
    @{Row=51; E="  +2.34%  "}
)

# Columns are plain-text cells in the source sheet (inline strings), several of
# which look numeric ("312.84", "0.0910", ...). Force text via NumberFormat "@"
# before the assignment so Excel does not silently coerce them to numbers, then
# restore the default "Normal" style so no stray number formatting is left behind.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { Set-TextCell $ws.Cells.Item($u.Row, 2) $u.B }
    if ($u.ContainsKey("C")) { Set-TextCell $ws.Cells.Item($u.Row, 3) $u.C }
    if ($u.ContainsKey("D")) { Set-TextCell $ws.Cells.Item($u.Row, 4) $u.D }
    if ($u.ContainsKey("E")) { Set-TextCell $ws.Cells.Item($u.Row, 5) $u.E }
}
